$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.396.01'
$ws.Range('E2').Value = '  -4.48%  '
$ws.Range('D3').Value = '1.570.48'
$ws.Range('E3').Value = '  -4.63%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.20'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3680'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.96%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.63'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3374'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.176'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07580'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -6.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.09'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.056'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -5.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.868'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -6.70%  '
$ws.Range('E16').Value = '  -4.11%  '
$ws.Range('D17').Value = '1.568.39'
$ws.Range('E17').Value = '  -5.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.12'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -8.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06710'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.52%  '
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  -7.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.41'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5239'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -9.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.99'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('D25').Value = '22.405.12'
$ws.Range('E25').Value = '  -4.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.385'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.971'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.88'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '145.81'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.950'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '125.05'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.84%  '
$ws.Range('D32').Value = '1.746.30'
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.266'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -9.34%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.004'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.977'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.36'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -11.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08431'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02536'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2300'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.527'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06507'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.07%  '
$ws.Range('E42').Value = '  -9.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.248'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.99%  '
$ws.Range('E44').Value = '  -7.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.59'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -6.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9997'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6028'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -5.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.763'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.121'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '121.12'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.84%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.198'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.75%  '
